# Update cryptos list values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.681.73"
$ws.Range("E2").Value = "  -3.61%  "
$ws.Range("D3").Value = "1.743.03"
$ws.Range("E3").Value = "  -5.52%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "237.60"
$ws.Range("E5").Value = "  -8.68%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.4940"
$ws.Range("E7").Value = "  -6.37%  "
$ws.Range("D8").Value = "41.60"
$ws.Range("E8").Value = "  -7.65%  "
$ws.Range("D9").Value = "0.2404"
$ws.Range("E9").Value = "  -23.68%  "
$ws.Range("D10").Value = "0.05958"
$ws.Range("E10").Value = "  -12.44%  "
$ws.Range("D11").Value = "1.742.59"
$ws.Range("E11").Value = "  -5.64%  "
$ws.Range("D12").Value = "0.06846"
$ws.Range("E12").Value = "  -11.97%  "
$ws.Range("D13").Value = "14.64"
$ws.Range("E13").Value = "  -23.44%  "
$ws.Range("D14").Value = "4.459"
$ws.Range("E14").Value = "  -11.17%  "
$ws.Range("D15").Value = "77.13"
$ws.Range("E15").Value = "  -12.67%  "
$ws.Range("D16").Value = "0.5821"
$ws.Range("E16").Value = "  -25.90%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "25.723.55"
$ws.Range("E19").Value = "  -3.55%  "
$ws.Range("D20").Value = "11.46"
$ws.Range("E20").Value = "  -17.65%  "
$ws.Range("D21").Value = "0.000006446"
$ws.Range("E21").Value = "  -18.70%  "
$ws.Range("D22").Value = "1.961.50"
$ws.Range("E22").Value = "  -6.11%  "
$ws.Range("D23").Value = "3.955"
$ws.Range("E23").Value = "  -14.23%  "
$ws.Range("D24").Value = "5.019"
$ws.Range("E24").Value = "  -16.29%  "
$ws.Range("D25").Value = "7.793"
$ws.Range("E25").Value = "  -16.74%  "
$ws.Range("D26").Value = "136.45"
$ws.Range("E26").Value = "  -4.57%  "
$ws.Range("E27").Value = "  -12.59%  "
$ws.Range("D28").Value = "1.831"
$ws.Range("E28").Value = "  -17.78%  "
$ws.Range("D29").Value = "14.51"
$ws.Range("E29").Value = "  -14.90%  "
$ws.Range("D30").Value = "100.70"
$ws.Range("E30").Value = "  -9.29%  "
$ws.Range("D31").Value = "3.790"
$ws.Range("E31").Value = "  -9.99%  "
$ws.Range("D32").Value = "0.08118"
$ws.Range("D33").Value = "3.346"
$ws.Range("E33").Value = "  -18.19%  "
$ws.Range("D34").Value = "0.04372"
$ws.Range("E34").Value = "  -10.52%  "
$ws.Range("D35").Value = "0.9994"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "2.639"
$ws.Range("E36").Value = "  -7.85%  "
$ws.Range("D37").Value = "1.022"
$ws.Range("E37").Value = "  -10.58%  "
$ws.Range("D38").Value = "0.6061"
$ws.Range("E38").Value = "  -17.25%  "
$ws.Range("D39").Value = "2.693"
$ws.Range("E39").Value = "  -13.39%  "
$ws.Range("D40").Value = "2.077"
$ws.Range("E40").Value = "  -10.22%  "
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "102.93"
$ws.Range("E42").Value = "  -6.27%  "
$ws.Range("E43").Value = "  -14.42%  "
$ws.Range("D44").Value = "0.7788"
$ws.Range("E44").Value = "  -14.12%  "
$ws.Range("D45").Value = "5.127"
$ws.Range("E45").Value = "  -13.74%  "
$ws.Range("E46").Value = "  -21.97%  "
$ws.Range("D47").Value = "0.05107"
$ws.Range("E47").Value = "  -12.44%  "
$ws.Range("D48").Value = "5.978"
$ws.Range("E48").Value = "  -22.60%  "
$ws.Range("D49").Value = "0.1064"
$ws.Range("E49").Value = "  -14.62%  "
$ws.Range("D50").Value = "30.01"
$ws.Range("E50").Value = "  -13.86%  "
$ws.Range("D51").Value = "52.66"
$ws.Range("E51").Value = "  -12.28%  "
